$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so Excel does not coerce
# numeric-looking strings (e.g. "1.00", "0.0750") into floating point numbers.
$touched = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","E7","E8","D9","E9","E10","E11","D12","E12","E13","D14","E14","D15","E15","E16","D17","E17","E18","D19","D20","E20","D21","E21","E22","D23","E23","E24","D25","E25","E26","E27","E28","E29","B30","C30","D30","E30","B31","C31","D31","E31","D32","E32","D33","E33","E34","E35","D36","E36","D37","E37","E38","D39","E39","D41","E41","E42","E43","E44","E45","D46","E46","E47","E48","D49","E49","E50","D51","E51")
foreach ($addr in $touched) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '60.786.67'
$ws.Range('E2').Value = '  -3.02%  '
$ws.Range('D3').Value = '3.346.72'
$ws.Range('E3').Value = '  -2.79%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '566.60'
$ws.Range('E5').Value = '  -2.24%  '
$ws.Range('D6').Value = '146.31'
$ws.Range('E6').Value = '  -0.89%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.47%  '
$ws.Range('D9').Value = '7.92'
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('E10').Value = '  -1.37%  '
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('D12').Value = '3.916.97'
$ws.Range('E12').Value = '  -2.84%  '
$ws.Range('E13').Value = '  +1.18%  '
$ws.Range('D14').Value = '27.61'
$ws.Range('E14').Value = '  -2.12%  '
$ws.Range('D15').Value = '3.349.18'
$ws.Range('E15').Value = '  -2.86%  '
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').Value = '60.788.72'
$ws.Range('E17').Value = '  -3.07%  '
$ws.Range('E18').Value = '  -1.44%  '
$ws.Range('D19').Value = '14.47'
$ws.Range('D20').Value = '8.88'
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('D21').Value = '375.49'
$ws.Range('E21').Value = '  -2.90%  '
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').Value = '74.79'
$ws.Range('E23').Value = '  -0.64%  '
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = '3.497.96'
$ws.Range('E25').Value = '  -2.33%  '
$ws.Range('E26').Value = '  -5.76%  '
$ws.Range('E27').Value = '  -3.99%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('E29').Value = '  -3.81%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '2.09'
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').Value = '7.72'
$ws.Range('E32').Value = '  -3.18%  '
$ws.Range('D33').Value = '22.78'
$ws.Range('E33').Value = '  -1.80%  '
$ws.Range('E34').Value = '  -3.62%  '
$ws.Range('E35').Value = '  -0.34%  '
$ws.Range('D36').Value = '168.63'
$ws.Range('E36').Value = '  -0.52%  '
$ws.Range('D37').Value = '1.55'
$ws.Range('E37').Value = '  -4.78%  '
$ws.Range('E38').Value = '  -2.32%  '
$ws.Range('D39').Value = '27.95'
$ws.Range('E39').Value = '  -12.26%  '
$ws.Range('D41').Value = '0.0750'
$ws.Range('E41').Value = '  -2.82%  '
$ws.Range('E42').Value = '  -3.45%  '
$ws.Range('E43').Value = '  -1.50%  '
$ws.Range('E44').Value = '  -5.21%  '
$ws.Range('E45').Value = '  -3.73%  '
$ws.Range('D46').Value = '2.459.15'
$ws.Range('E46').Value = '  -4.14%  '
$ws.Range('E47').Value = '  -4.09%  '
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('D49').Value = '22.22'
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('E50').Value = '  -2.04%  '
$ws.Range('D51').Value = '0.816'
$ws.Range('E51').Value = '  +0.58%  '
